$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translations")
$ws.Rows.Item(707).Resize(2).Delete()
